# Update "want-to-go" head-count figures (column F) on the "展览" (sheet 1),
# "演出" (sheet 2) and "全部类型" (sheet 4) worksheets, per the regenerated
# gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value  = 2900
$ws1.Range("F6").Value  = 68
$ws1.Range("F8").Value  = 2260
$ws1.Range("F9").Value  = 1481
$ws1.Range("F10").Value = 41
$ws1.Range("F11").Value = 540
$ws1.Range("F13").Value = 2595
$ws1.Range("F15").Value = 1430
$ws1.Range("F16").Value = 5571
$ws1.Range("F18").Value = 5453
$ws1.Range("F19").Value = 2047
$ws1.Range("F20").Value = 2959
$ws1.Range("F21").Value = 3389
$ws1.Range("F23").Value = 1655
$ws1.Range("F24").Value = 33
$ws1.Range("F25").Value = 277
$ws1.Range("F26").Value = 853
$ws1.Range("F27").Value = 151
$ws1.Range("F28").Value = 12
$ws1.Range("F29").Value = 343
$ws1.Range("F31").Value = 2229
$ws1.Range("F33").Value = 135
$ws1.Range("F34").Value = 316
$ws1.Range("F35").Value = 834
$ws1.Range("F36").Value = 171
$ws1.Range("F37").Value = 401
$ws1.Range("F38").Value = 467

# --- Sheet 2: 演出 (Performance) ---------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value  = 62
$ws2.Range("F20").Value = 6

# --- Sheet 4: 全部类型 (All types) --------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value  = 62
$ws4.Range("F10").Value = 2900
$ws4.Range("F11").Value = 68
$ws4.Range("F12").Value = 2260
$ws4.Range("F13").Value = 1481
$ws4.Range("F14").Value = 41
$ws4.Range("F15").Value = 540
$ws4.Range("F18").Value = 2595
$ws4.Range("F19").Value = 1430
$ws4.Range("F24").Value = 5572
$ws4.Range("F26").Value = 5454
$ws4.Range("F27").Value = 2047
$ws4.Range("F28").Value = 2959
$ws4.Range("F29").Value = 3389
$ws4.Range("F34").Value = 1655
$ws4.Range("F36").Value = 277
$ws4.Range("F37").Value = 853
$ws4.Range("F38").Value = 151
$ws4.Range("F39").Value = 12
$ws4.Range("F40").Value = 343
$ws4.Range("F42").Value = 2229
$ws4.Range("F44").Value = 135
$ws4.Range("F45").Value = 316
$ws4.Range("F46").Value = 834
$ws4.Range("F47").Value = 171
$ws4.Range("F48").Value = 401
$ws4.Range("F49").Value = 467
